{"js": "// Fix the misspelling \"arraign\u00e9\" -> \"araign\u00e9e\" and drop the now-stale\n// spell-check bookmarks (w:proofErr spellStart/spellEnd) that used to\n// flag the typo, mirroring the author's OOXML diff.\nconst misspelled = \"arraign\u00e9\";\nconst corrected = \"araign\u00e9e\";\n\nconst results = context.document.body.search(misspelled, {\n  matchCase: true,\n  matchWholeWord: true\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Text \"${misspelled}\" not found in document body.`);\n}\n\n// Work at the paragraph level: re-serializing the paragraph's OOXML\n// through getOoxml()/insertOoxml() naturally drops the stale\n// <w:proofErr/> spell-check markers that bracket the misspelled run,\n// while every other run/paragraph property (rsids, tabs, etc.) is\n// preserved untouched.\nconst paragraph = results.items[0].paragraphs.getFirst();\nconst ooxml = paragraph.getOoxml();\nawait context.sync();\n\nlet xml = ooxml.value;\nxml = xml.replace(`>${misspelled}<`, `>${corrected}<`);\n// Strip the paraId/textId stamps the OOXML export synthesizes for the\n// fragment - they are not present on the original paragraph.\nxml = xml.replace(/\\s+w14:paraId=\"[^\"]*\"/, \"\").replace(/\\s+w14:textId=\"[^\"]*\"/, \"\");\n\nparagraph.insertOoxml(xml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fix the misspelling \"arraign\u00e9\" -> \"araign\u00e9e\" and drop the now-stale\n# spell-check bookmarks (w:proofErr spellStart/spellEnd) that used to\n# flag the typo, mirroring the author's OOXML diff.\n$d = $word.ActiveDocument\n$misspelled = \"arraign\u00e9\"\n$corrected = \"araign\u00e9e\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains($misspelled)) {\n        $target = $p.Range\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Paragraph containing '$misspelled' not found\"\n}\n\n# Round-trip the paragraph through WordOpenXML (read-only) / InsertXML\n# (write): re-serializing the paragraph drops the stale <w:proofErr/>\n# spell-check markers that bracket the misspelled run, while every\n# other run/paragraph property (rsids, tabs, etc.) is preserved as-is.\n$xml = $target.WordOpenXML\n$xml = $xml.Replace(\">$misspelled<\", \">$corrected<\")\n$xml = $xml -replace ' w14:paraId=\"[^\"]*\"', ''\n$xml = $xml -replace ' w14:textId=\"[^\"]*\"', ''\n\n$target.InsertXML($xml)\n"}
